# Commit: Tue, Jul 28, 2020 8:07:42 PM
#
# Change the table style applied to the table on slide 5
# ("B1- TYPES OF FINANCIAL DOCUMENTS") from the default table style
# to a different built-in table style (identified by its style GUID).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

foreach ($shp in $s.Shapes) {
    if ($shp.HasTable) {
        $tbl = $shp.Table
        $tbl.ApplyStyle("{55CD6981-3575-4CB7-954C-FD5E6B28472B}")
    }
}
